# Add a new "2021" column (K) to the table, mirroring the formatting of the
# existing "2020" column (J), and fill in the reported values for that year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/borders) of column J (rows 3-10, the table
# body+header+footer rows) into the new column K.
$ws.Range("J3:J10").Copy()
$ws.Range("K3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new column's values.
$ws.Range("K4").Value = 2021
$ws.Range("K5").Value = 375
$ws.Range("K6").Value = "-"
$ws.Range("K7").Value = 5
$ws.Range("K8").Value = "-"
$ws.Range("K9").Value = 18
$ws.Range("K10").Value = 150

# Match the saved selection state (active cell K7) seen in the workbook.
$ws.Range("K7").Select()
